$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the string "Teh ..." labels in column A (rows 2-21) with plain numbers.
$values = @{
    2  = 3
    3  = 4
    4  = 3
    5  = 9
    6  = 2
    7  = 1
    8  = 9
    9  = 5
    10 = 6
    11 = 7
    12 = 9
    13 = 8
    14 = 2
    15 = 2
    16 = 4
    17 = 4
    18 = 7
    19 = 8
    20 = 6
    21 = 4
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}

# Update the sheet view / selection to match the saved state.
$excel.ActiveWindow.ScrollRow = 12
$ws.Range("A21").Select()
